$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Row 2 ---
$ws1.Range("B2").Value = "Planche à Voile"
$ws1.Range("G2").Value = "Yzc2MDU5NmU4N2MyZDdkNDIwZjIwNz"

# --- Row 3 ---
$ws1.Range("A3").Value = "M"
$ws1.Range("B3").Value = "Paquebot"
$ws1.Range("G3").Value = "Yzc2MDU5NmU4N2MyZDdkNDIwZjIwNz"

# --- Row 4 ---
$ws1.Range("A4").Value = "M"
$ws1.Range("B4").Value = "Voilier"
$ws1.Range("D4").Value = 35351
$ws1.Range("G4").Value = "Yzc2MDU5NmU4N2MyZDdkNDIwZjIwNz"

# --- Hyperlinks on the Email column: rebuild from scratch ---
$ws1.Hyperlinks.Delete()

$ws1.Hyperlinks.Add($ws1.Range("E2"), "mailto:bb@bbmlm.ci19")
$ws1.Hyperlinks.Add($ws1.Range("E3"), "mailto:bb@bbmlm.ci19", "", "", "bb@bbmlm.ci19")
$ws1.Hyperlinks.Add($ws1.Range("E4"), "mailto:bb@bbmlm.ci19", "", "", "bb@bbmlm.ci19")

$ws1.Range("E2").Value = "bb@bbmlm.ci19"
$ws1.Range("E3").Value = "bb@bbmlm.ci20"
$ws1.Range("E4").Value = "bb@bbmlm.ci21"

# --- Remove now-unused trailing blank rows (14:24), keep 5:13 ---
$ws1.Rows("14:24").Delete()

# --- Update the remembered selection to match the post-edit state ---
$ws1.Range("G7").Select() | Out-Null
